$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4263.6665
$ws.Range("J40").Value = 4874.5
$ws.Range("L40").Value = 4874.5
$ws.Range("N40").Value = -5224.5
$ws.Range("H70").Value = 2960.25
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 2960.25
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H76").Value = 3059
$ws.Range("I76").Value = 2850
$ws.Range("K76").Value = 2850
$ws.Range("M76").Value = -2535
$ws.Range("H79").Value = 3059
$ws.Range("I79").Value = 2850
$ws.Range("K79").Value = 2850
$ws.Range("M79").Value = -1758
$ws.Range("H98").Value = 1614.2903
$ws.Range("I98").Value = 728.0741
$ws.Range("K98").Value = 728.0741
$ws.Range("M98").Value = 769.9259
$ws.Range("H106").Value = 90911130
$ws.Range("I106").Value = 100001890
$ws.Range("K106").Value = 100001890
$ws.Range("M106").Value = -100001259
$ws.Range("H107").Value = 435.3125
$ws.Range("I107").Value = 438.33334
$ws.Range("K107").Value = 438.33334
$ws.Range("M107").Value = 1481.66666
$ws.Range("H122").Value = 1614.2903
$ws.Range("I122").Value = 728.0741
$ws.Range("K122").Value = 2184.2223
$ws.Range("M122").Value = 265.7776999999996
$ws.Range("H135").Value = 1271.4755
$ws.Range("I135").Value = 1030.3773
$ws.Range("J135").Value = 2868.75
$ws.Range("K135").Value = 9273.395700000001
$ws.Range("L135").Value = 25818.75
$ws.Range("M135").Value = -6738.395700000001
$ws.Range("N135").Value = -30888.75
$ws.Range("H138").Value = 5931.74
$ws.Range("J138").Value = 3827.6553
$ws.Range("L138").Value = 11482.9659
$ws.Range("N138").Value = -21762.9659
$ws.Range("H139").Value = 153993.8
$ws.Range("J139").Value = 153993.8
$ws.Range("L139").Value = 153993.8
$ws.Range("N139").Value = -164273.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1240.0212
$ws.Range("I2").Value = 912.65717
$ws.Range("K2").Value = 912.65717
$ws.Range("M2").Value = -799.65717
$ws.Range("H25").Value = 12959
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 14785.5
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 14785.5
$ws.Range("M25").Value = -1598
$ws.Range("N25").Value = -15589.5
$ws.Range("H32").Value = 139653.05
$ws.Range("I32").Value = 147492.7
$ws.Range("K32").Value = 147492.7
$ws.Range("M32").Value = -147205.7
$ws.Range("H45").Value = 30265.666
$ws.Range("I45").Value = 45015.87
$ws.Range("K45").Value = 45015.87
$ws.Range("M45").Value = -44638.87
$ws.Range("H61").Value = 4916277.5
$ws.Range("I61").Value = 27372.8
$ws.Range("K61").Value = 27372.8
$ws.Range("M61").Value = -27160.8
$ws.Range("H74").Value = 1018101.56
$ws.Range("I74").Value = 8564.200000000001
$ws.Range("J74").Value = 1859382.8
$ws.Range("K74").Value = 8564.200000000001
$ws.Range("L74").Value = 1859382.8
$ws.Range("M74").Value = -7690.200000000001
$ws.Range("N74").Value = -1861130.8
$ws.Range("H77").Value = 1018101.56
$ws.Range("I77").Value = 8564.200000000001
$ws.Range("J77").Value = 1859382.8
$ws.Range("K77").Value = 42821
$ws.Range("L77").Value = 9296914
$ws.Range("M77").Value = -38453
$ws.Range("N77").Value = -9305650
$ws.Range("H116").Value = 1240.0212
$ws.Range("I116").Value = 912.65717
$ws.Range("K116").Value = 912.65717
$ws.Range("M116").Value = 1381.34283
$ws.Range("H136").Value = 4916277.5
$ws.Range("I136").Value = 27372.8
$ws.Range("K136").Value = 82118.39999999999
$ws.Range("M136").Value = -79568.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1240.0212
$ws.Range("I3").Value = 912.65717
$ws.Range("K3").Value = 912.65717
$ws.Range("M3").Value = -798.65717
$ws.Range("H20").Value = 1060.725
$ws.Range("I20").Value = 976.2857
$ws.Range("J20").Value = 1257.75
$ws.Range("K20").Value = 976.2857
$ws.Range("L20").Value = 1257.75
$ws.Range("M20").Value = -729.2857
$ws.Range("N20").Value = -1751.75
$ws.Range("H107").Value = 10391.939
$ws.Range("I107").Value = 11326.786
$ws.Range("K107").Value = 11326.786
$ws.Range("M107").Value = -9406.786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2195.5881
$ws.Range("I16").Value = 2110.0833
$ws.Range("K16").Value = 2110.0833
$ws.Range("M16").Value = -1823.0833
$ws.Range("H58").Value = 3110.5454
$ws.Range("I58").Value = 1949.75
$ws.Range("K58").Value = 1949.75
$ws.Range("M58").Value = -1746.75
$ws.Range("H62").Value = 8926.444
$ws.Range("J62").Value = 6001
$ws.Range("L62").Value = 6001
$ws.Range("N62").Value = -7249
$ws.Range("H65").Value = 8926.444
$ws.Range("J65").Value = 6001
$ws.Range("L65").Value = 30005
$ws.Range("N65").Value = -36245
$ws.Range("H68").Value = 79499.5
$ws.Range("J68").Value = 69999
$ws.Range("L68").Value = 69999
$ws.Range("N68").Value = -71497
$ws.Range("H71").Value = 79499.5
$ws.Range("J71").Value = 69999
$ws.Range("L71").Value = 209997
$ws.Range("N71").Value = -217485
$ws.Range("H113").Value = 2195.5881
$ws.Range("I113").Value = 2110.0833
$ws.Range("K113").Value = 2110.0833
$ws.Range("M113").Value = 59.91670000000022
$ws.Range("H134").Value = 6289.615
$ws.Range("I134").Value = 6806.6
$ws.Range("K134").Value = 20419.8
$ws.Range("M134").Value = -17884.8
$ws.Range("H136").Value = 3110.5454
$ws.Range("I136").Value = 1949.75
$ws.Range("K136").Value = 5849.25
$ws.Range("M136").Value = -3299.25
$ws.Range("H137").Value = 125294
$ws.Range("J137").Value = 125294
$ws.Range("L137").Value = 125294
$ws.Range("N137").Value = -135494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 361.5
$ws.Range("I8").Value = 361.5
$ws.Range("K8").Value = 1084.5
$ws.Range("M8").Value = -945.5
$ws.Range("H131").Value = 4063793.8
$ws.Range("I131").Value = 5684044.5
$ws.Range("J131").Value = 3301322.5
$ws.Range("K131").Value = 17052133.5
$ws.Range("L131").Value = 9903967.5
$ws.Range("M131").Value = -17047093.5
$ws.Range("N131").Value = -9914047.5
$ws.Range("H132").Value = 2429.6956
$ws.Range("I132").Value = 2273
$ws.Range("K132").Value = 20457
$ws.Range("M132").Value = -17927

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6145425.5
$ws.Range("I80").Value = 128661.734
$ws.Range("J80").Value = 22828272
$ws.Range("K80").Value = 128661.734
$ws.Range("L80").Value = 22828272
$ws.Range("M80").Value = -127663.734
$ws.Range("N80").Value = -22830268
$ws.Range("H83").Value = 6145425.5
$ws.Range("I83").Value = 128661.734
$ws.Range("J83").Value = 22828272
$ws.Range("K83").Value = 643308.6699999999
$ws.Range("L83").Value = 114141360
$ws.Range("M83").Value = -638316.6699999999
$ws.Range("N83").Value = -114151344
$ws.Range("H95").Value = 39997.5
$ws.Range("J95").Value = 39997.5
$ws.Range("L95").Value = 39997.5
$ws.Range("N95").Value = -45489.5
$ws.Range("H97").Value = 605.3226
$ws.Range("I97").Value = 638.6842
$ws.Range("J97").Value = 552.5
$ws.Range("K97").Value = 638.6842
$ws.Range("L97").Value = 552.5
$ws.Range("M97").Value = -142.6842
$ws.Range("N97").Value = -1544.5
$ws.Range("H102").Value = 20001708
$ws.Range("I102").Value = 23811006
$ws.Range("K102").Value = 23811006
$ws.Range("M102").Value = -23809384
$ws.Range("H113").Value = 2500.8918
$ws.Range("I113").Value = 1853.4231
$ws.Range("K113").Value = 1853.4231
$ws.Range("M113").Value = 316.5769
$ws.Range("H123").Value = 26997.8
$ws.Range("J123").Value = 26997.8
$ws.Range("L123").Value = 26997.8
$ws.Range("N123").Value = -31897.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15544.211
$ws.Range("I7").Value = 15058.097
$ws.Range("K7").Value = 15058.097
$ws.Range("M7").Value = -14946.097
$ws.Range("H40").Value = 4418.8
$ws.Range("I40").Value = 4285.952
$ws.Range("K40").Value = 4285.952
$ws.Range("M40").Value = -4149.952
$ws.Range("H46").Value = 6483
$ws.Range("I46").Value = 20453.6
$ws.Range("J46").Value = 1493.5
$ws.Range("K46").Value = 20453.6
$ws.Range("L46").Value = 1493.5
$ws.Range("M46").Value = -20265.6
$ws.Range("N46").Value = -1869.5
$ws.Range("H93").Value = 1997.2727
$ws.Range("I93").Value = 2415.75
$ws.Range("J93").Value = 1495.1
$ws.Range("K93").Value = 2415.75
$ws.Range("L93").Value = 1495.1
$ws.Range("M93").Value = -1167.75
$ws.Range("N93").Value = -3991.1
$ws.Range("H100").Value = 1463.619
$ws.Range("I100").Value = 1458.909
$ws.Range("J100").Value = 1468.8
$ws.Range("K100").Value = 1458.909
$ws.Range("L100").Value = 1468.8
$ws.Range("M100").Value = -917.9090000000001
$ws.Range("N100").Value = -2550.8
$ws.Range("H126").Value = 15544.211
$ws.Range("I126").Value = 15058.097
$ws.Range("K126").Value = 45174.291
$ws.Range("M126").Value = -42704.291
$ws.Range("H135").Value = 55950.08
$ws.Range("J135").Value = 55950.08
$ws.Range("L135").Value = 55950.08
$ws.Range("N135").Value = -66090.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 37341
$ws.Range("J54").Value = 39062.312
$ws.Range("L54").Value = 39062.312
$ws.Range("N54").Value = -40102.312
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 50000
$ws.Range("K60").Value = 50000
$ws.Range("M60").Value = -49178
$ws.Range("H81").Value = 1500.6666
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 1502
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 3004
$ws.Range("M81").Value = -1939
$ws.Range("N81").Value = -5126
$ws.Range("H84").Value = 1500.6666
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 1502
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 15020
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -25628
$ws.Range("H100").Value = 3617
$ws.Range("I100").Value = 3617
$ws.Range("K100").Value = 7234
$ws.Range("M100").Value = -6693
$ws.Range("H107").Value = 142858260
$ws.Range("I107").Value = 1428.2
$ws.Range("K107").Value = 4284.6
$ws.Range("M107").Value = -2364.6
$ws.Range("H122").Value = 2427.9
$ws.Range("I122").Value = 2380.24
$ws.Range("J122").Value = 2666.2
$ws.Range("K122").Value = 7140.719999999999
$ws.Range("L122").Value = 7998.599999999999
$ws.Range("M122").Value = -4690.719999999999
$ws.Range("N122").Value = -12898.6
$ws.Range("H132").Value = 29483
$ws.Range("I132").Value = 48084.668
$ws.Range("K132").Value = 144254.004
$ws.Range("M132").Value = -141724.004
$ws.Range("H136").Value = 5855.6
$ws.Range("I136").Value = 6825.2
$ws.Range("K136").Value = 20475.6
$ws.Range("M136").Value = -17925.6
